# Update the recorded bookmark PDF path (shared string in A1) and move the
# sheet's recorded selection back onto A1 (closest reproducible approximation
# of the author's cleaned-up view state - the source file had a stray
# "B3" selection left over from editing).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "C:\Users\martin.asenov\Desktop\Scanned\A7-MHH 106-1 to 106-39.pdf"
[void]$ws.Range("A1").Select()
